$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New day's data (2025-09-21, serial date 45921) for both charging stations
$ws.Range("A42").Value = 45921
$ws.Range("B42").Value = "四方坪站"
$ws.Range("C42").Value = 8895.42
$ws.Range("D42").Value = 7042.34
$ws.Range("E42").Value = 3013.42
$ws.Range("F42").Value = 370

$ws.Range("A43").Value = 45921
$ws.Range("B43").Value = "高岭站"
$ws.Range("C43").Value = 4272.43
$ws.Range("D43").Value = 3339.78
$ws.Range("E43").Value = 1084.7
$ws.Range("F43").Value = 147

# Match the number formats used by the existing rows in each column
# (xlPasteFormats copies just the cell formatting, not the values we just set)
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$ws.Range("A41").Copy()
$ws.Range("A42:A43").PasteSpecial($xlPasteFormats)

$ws.Range("C41:E41").Copy()
$ws.Range("C42:E43").PasteSpecial($xlPasteFormats)

$ws.Range("F41").Copy()
$ws.Range("F42:F43").PasteSpecial($xlPasteFormats)

$ws.Range("H41").Select()
